$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 800.63635
$ws.Range("I6").Value = 576.25
$ws.Range("K6").Value = 1728.75
$ws.Range("M6").Value = -1616.75
$ws.Range("H9").Value = 181.05882
$ws.Range("J9").Value = 450
$ws.Range("L9").Value = 450
$ws.Range("N9").Value = -788
$ws.Range("H12").Value = 1317.875
$ws.Range("J12").Value = 1579.8
$ws.Range("L12").Value = 1579.8
$ws.Range("N12").Value = -1919.8
$ws.Range("H21").Value = 25000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30936
$ws.Range("H23").Value = 25000
$ws.Range("J23").Value = 30000
$ws.Range("L23").Value = 30000
$ws.Range("N23").Value = -30468
$ws.Range("H29").Value = 4000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 7964.75
$ws.Range("I62").Value = 7987.3335
$ws.Range("J62").Value = 7897
$ws.Range("K62").Value = 7987.3335
$ws.Range("L62").Value = 7897
$ws.Range("M62").Value = -7363.3335
$ws.Range("N62").Value = -9145
$ws.Range("H65").Value = 7964.75
$ws.Range("I65").Value = 7987.3335
$ws.Range("J65").Value = 7897
$ws.Range("K65").Value = 39936.6675
$ws.Range("L65").Value = 39485
$ws.Range("M65").Value = -36816.6675
$ws.Range("N65").Value = -45725
$ws.Range("H86").Value = 5271.2856
$ws.Range("I86").Value = 4939.8
$ws.Range("K86").Value = 4939.8
$ws.Range("M86").Value = -3816.8
$ws.Range("H89").Value = 5271.2856
$ws.Range("I89").Value = 4939.8
$ws.Range("K89").Value = 24699
$ws.Range("M89").Value = -19083
$ws.Range("H95").Value = 43428.57
$ws.Range("J95").Value = 43428.57
$ws.Range("L95").Value = 43428.57
$ws.Range("N95").Value = -48920.57
$ws.Range("H137").Value = 2091.3684
$ws.Range("J137").Value = 2411.6155
$ws.Range("L137").Value = 7234.8465
$ws.Range("N137").Value = -12334.8465
$ws.Range("H138").Value = 2354.6829
$ws.Range("I138").Value = 1802.5217
$ws.Range("J138").Value = 3060.2222
$ws.Range("K138").Value = 5407.5651
$ws.Range("L138").Value = 9180.6666
$ws.Range("M138").Value = -267.5650999999998
$ws.Range("N138").Value = -19460.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2401.9395
$ws.Range("I32").Value = 2524.4727
$ws.Range("K32").Value = 2524.4727
$ws.Range("M32").Value = -2237.4727
$ws.Range("H61").Value = 3121.3901
$ws.Range("I61").Value = 2765.1936
$ws.Range("K61").Value = 2765.1936
$ws.Range("M61").Value = -2553.1936
$ws.Range("H101").Value = 69601.664
$ws.Range("J101").Value = 69601.664
$ws.Range("L101").Value = 69601.664
$ws.Range("N101").Value = -76091.664
$ws.Range("H132").Value = 4234.919
$ws.Range("I132").Value = 2948.7307
$ws.Range("K132").Value = 8846.1921
$ws.Range("M132").Value = -6316.1921
$ws.Range("H136").Value = 3121.3901
$ws.Range("I136").Value = 2765.1936
$ws.Range("K136").Value = 8295.5808
$ws.Range("M136").Value = -5745.5808

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3809.6545
$ws.Range("I134").Value = 3653.795
$ws.Range("K134").Value = 10961.385
$ws.Range("M134").Value = -8426.385

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3687.2917
$ws.Range("J31").Value = 5803.4165
$ws.Range("L31").Value = 5803.4165
$ws.Range("N31").Value = -6393.4165
$ws.Range("H34").Value = 3687.2917
$ws.Range("J34").Value = 5803.4165
$ws.Range("L34").Value = 5803.4165
$ws.Range("N34").Value = -6207.4165
$ws.Range("H58").Value = 2534.6562
$ws.Range("I58").Value = 2782.8333
$ws.Range("K58").Value = 2782.8333
$ws.Range("M58").Value = -2579.8333
$ws.Range("H136").Value = 2534.6562
$ws.Range("I136").Value = 2782.8333
$ws.Range("K136").Value = 8348.499899999999
$ws.Range("M136").Value = -5798.499899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4771.3887
$ws.Range("I14").Value = 4771.3887
$ws.Range("K14").Value = 14314.1661
$ws.Range("M14").Value = -14141.1661
$ws.Range("H37").Value = 129313.5
$ws.Range("J37").Value = 129313.5
$ws.Range("L37").Value = 387940.5
$ws.Range("N37").Value = -388164.5
$ws.Range("H86").Value = 569.17645
$ws.Range("I86").Value = 358.7143
$ws.Range("J86").Value = 1551.3334
$ws.Range("K86").Value = 1076.1429
$ws.Range("L86").Value = 4654.0002
$ws.Range("M86").Value = 109.8571000000002
$ws.Range("N86").Value = -7026.0002
$ws.Range("H89").Value = 569.17645
$ws.Range("I89").Value = 358.7143
$ws.Range("J89").Value = 1551.3334
$ws.Range("K89").Value = 3228.4287
$ws.Range("L89").Value = 13962.0006
$ws.Range("M89").Value = 2699.5713
$ws.Range("N89").Value = -25818.0006
$ws.Range("H131").Value = 1924.619
$ws.Range("J131").Value = 2982.889
$ws.Range("L131").Value = 8948.667000000001
$ws.Range("N131").Value = -19028.667
$ws.Range("H136").Value = 3754.2
$ws.Range("I136").Value = 3657.5454
$ws.Range("J136").Value = 4020
$ws.Range("K136").Value = 10972.6362
$ws.Range("L136").Value = 12060
$ws.Range("M136").Value = -5872.636200000001
$ws.Range("N136").Value = -22260

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3364.9092
$ws.Range("I126").Value = 3401.4
$ws.Range("K126").Value = 10204.2
$ws.Range("M126").Value = -7734.200000000001
$ws.Range("H132").Value = 3144.0667
$ws.Range("I132").Value = 3284.389
$ws.Range("J132").Value = 2933.5833
$ws.Range("K132").Value = 9853.167000000001
$ws.Range("L132").Value = 8800.749899999999
$ws.Range("M132").Value = -7323.167000000001
$ws.Range("N132").Value = -13860.7499

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2557.4285
$ws.Range("I22").Value = 1299
$ws.Range("K22").Value = 1299
$ws.Range("M22").Value = -1004
$ws.Range("H27").Value = 2557.4285
$ws.Range("I27").Value = 1299
$ws.Range("K27").Value = 1299
$ws.Range("M27").Value = -1192
$ws.Range("H46").Value = 4460
$ws.Range("I46").Value = 1375
$ws.Range("J46").Value = 7545
$ws.Range("K46").Value = 1375
$ws.Range("L46").Value = 7545
$ws.Range("M46").Value = -1187
$ws.Range("N46").Value = -7921
$ws.Range("H136").Value = 10757421
$ws.Range("I136").Value = 3300.4211
$ws.Range("K136").Value = 9901.263300000001
$ws.Range("M136").Value = -7351.263300000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 746.35297
$ws.Range("I113").Value = 803.4286
$ws.Range("J113").Value = 480
$ws.Range("K113").Value = 2410.2858
$ws.Range("L113").Value = 1440
$ws.Range("M113").Value = -240.2857999999997
$ws.Range("N113").Value = -5780
$ws.Range("H136").Value = 3815.9048
$ws.Range("I136").Value = 4193
$ws.Range("K136").Value = 12579
$ws.Range("M136").Value = -10029
